$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.163.26'
$ws.Range("E2").Value = '  -0.17%  '
$ws.Range("D3").Value = '2.528.17'
$ws.Range("E3").Value = '  +0.29%  '
$ws.Range("D4").Value = '''0.998'
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").Value = '''535.84'
$ws.Range("E5").Value = '  -0.17%  '
$ws.Range("D6").Value = '''138.24'
$ws.Range("E6").Value = '  -1.23%  '
$ws.Range("D7").Value = '''0.998'
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("E8").Value = '  +0.95%  '
$ws.Range("D9").Value = '2.526.21'
$ws.Range("E9").Value = '  +0.05%  '
$ws.Range("E10").Value = '  +1.74%  '
$ws.Range("E11").Value = '  +0.03%  '
$ws.Range("E12").Value = '  -1.43%  '
$ws.Range("E13").Value = '  -2.63%  '
$ws.Range("D14").Value = '2.959.62'
$ws.Range("E14").Value = '  -0.25%  '
$ws.Range("D15").Value = '''23.30'
$ws.Range("E15").Value = '  +1.78%  '
$ws.Range("D16").Value = '59.131.92'
$ws.Range("E16").Value = '  -0.11%  '
$ws.Range("E17").Value = '  -0.43%  '
$ws.Range("D18").Value = '2.507.63'
$ws.Range("E18").Value = '  +0.01%  '
$ws.Range("D19").Value = '''11.11'
$ws.Range("E19").Value = '  +1.69%  '
$ws.Range("E20").Value = '  +0.71%  '
$ws.Range("D21").Value = '''325.12'
$ws.Range("E21").Value = '  +0.73%  '
$ws.Range("E22").Value = '  +0.44%  '
$ws.Range("D23").Value = '''5.90'
$ws.Range("E23").Value = '  +1.18%  '
$ws.Range("D24").Value = '''65.03'
$ws.Range("E24").Value = '  +4.35%  '
$ws.Range("E25").Value = '  -0.99%  '
$ws.Range("E26").Value = '  +0.63%  '
$ws.Range("E27").Value = '  +1.38%  '
$ws.Range("E28").Value = '  -2.00%  '
$ws.Range("D29").Value = '0.0₃0779'
$ws.Range("E29").Value = '  +1.62%  '
$ws.Range("D30").Value = '''6.73'
$ws.Range("E30").Value = '  -0.84%  '
$ws.Range("E31").Value = '  -1.83%  '
$ws.Range("D32").Value = '''168.23'
$ws.Range("E32").Value = '  +3.69%  '
$ws.Range("D33").Value = '''1.20'
$ws.Range("E33").Value = '  +6.39%  '
$ws.Range("E34").Value = '  -0.12%  '
$ws.Range("E35").Value = '  -2.34%  '
$ws.Range("D36").Value = '''18.60'
$ws.Range("E36").Value = '  +0.62%  '
$ws.Range("D37").Value = '''4.12'
$ws.Range("E37").Value = '  -2.28%  '
$ws.Range("E38").Value = '  -0.96%  '
$ws.Range("D39").Value = '''36.77'
$ws.Range("D40").Value = '''0.833'
$ws.Range("E40").Value = '  +3.58%  '
$ws.Range("E41").Value = '  +0.27%  '
$ws.Range("D42").Value = '''5.27'
$ws.Range("E42").Value = '  +0.50%  '
$ws.Range("D43").Value = '''282.53'
$ws.Range("E43").Value = '  -0.58%  '
$ws.Range("D44").Value = '''0.999'
$ws.Range("E44").Value = '  +0.17%  '
$ws.Range("D45").Value = '''130.91'
$ws.Range("E45").Value = '  +6.68%  '
$ws.Range("D46").Value = '''0.607'
$ws.Range("E46").Value = '  +1.71%  '
$ws.Range("D47").Value = '''10.87'
$ws.Range("E47").Value = '  +0.13%  '
$ws.Range("D49").Value = '''0.0514'
$ws.Range("E49").Value = '  +0.59%  '
$ws.Range("E50").Value = '  +0.37%  '
$ws.Range("D51").Value = '''17.42'
$ws.Range("E51").Value = '  -0.16%  '
